# remaining days column added -- adds the new payment record row (row 2)
# to Sheet1: date, member_id, member_name, package, amount, status.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date (text, e.g. "31-03-2025")
$ws.Range("A2").Value = "31-03-2025"
$ws.Range("A2").Style = "Normal"

# member_id -- stored as text ("1001"), not a number, in the source file.
# Format as text before assigning so Excel doesn't auto-coerce the
# numeric-looking string into a Number cell, then drop back to the
# default "Normal" style so no explicit style index is left on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1001"
$ws.Range("B2").Style = "Normal"

# member_name
$ws.Range("C2").Value = "Abdullah"
$ws.Range("C2").Style = "Normal"

# package
$ws.Range("D2").Value = "Gold"
$ws.Range("D2").Style = "Normal"

# amount -- numeric
$ws.Range("E2").Value = 50000
$ws.Range("E2").Style = "Normal"

# status
$ws.Range("F2").Value = "Paid"
$ws.Range("F2").Style = "Normal"
